$wb = $excel.ActiveWorkbook

# Sheet "Hoja1": update the conversion rate text in A1
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.35 = 12929.74 pesos`n✅ 12929.74 pesos = 3.33 = 943.04 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Sheet "tasas": update rate values
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 298.815
$ws2.Range("O10").Value = 3863.6
$ws2.Range("N12").Value = 3886.99
$ws2.Range("O12").Value = 283.5
